$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7666.3335
$ws.Range("I113").Value = 6726
$ws.Range("J113").Value = 7935
$ws.Range("K113").Value = 6726
$ws.Range("L113").Value = 7935
$ws.Range("M113").Value = -3472
$ws.Range("N113").Value = -14443

$ws.Range("H125").Value = 3366.4
$ws.Range("I125").Value = 4523.4287
$ws.Range("J125").Value = 666.6667
$ws.Range("K125").Value = 40710.85830000001
$ws.Range("L125").Value = 6000.0003
$ws.Range("M125").Value = -38250.85830000001
$ws.Range("N125").Value = -10920.0003

$ws.Range("H132").Value = 2832.4167
$ws.Range("I132").Value = 2999
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 8997
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -6467
$ws.Range("N132").Value = -8060

$ws.Range("H137").Value = 1544.7333
$ws.Range("I137").Value = 1102.2106
$ws.Range("J137").Value = 1868.1154
$ws.Range("K137").Value = 3306.6318
$ws.Range("L137").Value = 5604.3462
$ws.Range("M137").Value = -756.6318000000001
$ws.Range("N137").Value = -10704.3462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7404.7915
$ws.Range("I45").Value = 8846.684999999999
$ws.Range("J45").Value = 1925.6
$ws.Range("K45").Value = 8846.684999999999
$ws.Range("L45").Value = 1925.6
$ws.Range("M45").Value = -8469.684999999999
$ws.Range("N45").Value = -2679.6

$ws.Range("H106").Value = 41651.6
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 41651.6
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 41651.6
$ws.Range("N106").Value = -44175.6

$ws.Range("H109").Value = 28594.25
$ws.Range("I109").Value = 20000
$ws.Range("J109").Value = 31459
$ws.Range("K109").Value = 20000
$ws.Range("L109").Value = 31459
$ws.Range("M109").Value = -18613
$ws.Range("N109").Value = -34233

$ws.Range("H125").Value = 49579.4
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 49579.4
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 49579.4
$ws.Range("N125").Value = -59419.4

$ws.Range("H132").Value = 2623.5854
$ws.Range("I132").Value = 1401.875
$ws.Range("J132").Value = 6967.4443
$ws.Range("K132").Value = 4205.625
$ws.Range("L132").Value = 20902.3329
$ws.Range("M132").Value = -1675.625

$ws.Range("H139").Value = 38656
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 38656
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 38656
$ws.Range("N139").Value = -48936

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 22223276
$ws.Range("I20").Value = 32258956
$ws.Range("J20").Value = 1411.3572
$ws.Range("K20").Value = 32258956
$ws.Range("L20").Value = 1411.3572
$ws.Range("M20").Value = -32258709
$ws.Range("N20").Value = -1905.3572

$ws.Range("H59").Value = 50400
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 50400
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 50400
$ws.Range("N59").Value = -52094

$ws.Range("H94").Value = 1589.8148
$ws.Range("I94").Value = 1296.35
$ws.Range("J94").Value = 2428.2856
$ws.Range("K94").Value = 1296.35
$ws.Range("L94").Value = 2428.2856
$ws.Range("M94").Value = -845.3499999999999
$ws.Range("N94").Value = -3330.2856

$ws.Range("H105").Value = 14642.75
$ws.Range("I105").Value = 34985
$ws.Range("J105").Value = 2437.4
$ws.Range("K105").Value = 34985
$ws.Range("L105").Value = 2437.4
$ws.Range("M105").Value = -33238

$ws.Range("H134").Value = 10510.786
$ws.Range("I134").Value = 29380.75
$ws.Range("J134").Value = 2962.8
$ws.Range("K134").Value = 88142.25
$ws.Range("L134").Value = 8888.400000000001
$ws.Range("M134").Value = -85607.25
$ws.Range("N134").Value = -13958.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4151.9287
$ws.Range("I31").Value = 2265.7058
$ws.Range("J31").Value = 4756.9434
$ws.Range("K31").Value = 2265.7058
$ws.Range("L31").Value = 4756.9434
$ws.Range("M31").Value = -1970.7058
$ws.Range("N31").Value = -5346.9434

$ws.Range("H34").Value = 4151.9287
$ws.Range("I34").Value = 2265.7058
$ws.Range("J34").Value = 4756.9434
$ws.Range("K34").Value = 2265.7058
$ws.Range("L34").Value = 4756.9434
$ws.Range("M34").Value = -2063.7058
$ws.Range("N34").Value = -5160.9434

$ws.Range("H92").Value = 38000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 38000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 38000
$ws.Range("N92").Value = -42992

$ws.Range("H95").Value = 6434.4
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 6434.4
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 6434.4
$ws.Range("N95").Value = -11926.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 133423.34
$ws.Range("I5").Value = 12896.375
$ws.Range("J5").Value = 158797.45
$ws.Range("K5").Value = 38689.125
$ws.Range("L5").Value = 476392.35
$ws.Range("M5").Value = -38577.125
$ws.Range("N5").Value = -476616.35

$ws.Range("H35").Value = 563.3333
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 563.3333
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 1689.9999
$ws.Range("N35").Value = -2265.9999

$ws.Range("H68").Value = 3152.7273
$ws.Range("I68").Value = 4900.077
$ws.Range("J68").Value = 2016.95
$ws.Range("K68").Value = 14700.231
$ws.Range("L68").Value = 6050.85
$ws.Range("M68").Value = -13889.231
$ws.Range("N68").Value = -7672.85

$ws.Range("H71").Value = 3152.7273
$ws.Range("I71").Value = 4900.077
$ws.Range("J71").Value = 2016.95
$ws.Range("K71").Value = 44100.693
$ws.Range("L71").Value = 18152.55
$ws.Range("M71").Value = -40044.693
$ws.Range("N71").Value = -26264.55

$ws.Range("H131").Value = 13581491
$ws.Range("I131").Value = 8333894
$ws.Range("J131").Value = 14494117
$ws.Range("K131").Value = 25001682
$ws.Range("L131").Value = 43482351
$ws.Range("M131").Value = -24996642
$ws.Range("N131").Value = -43492431

$ws.Range("H135").Value = 133423.34
$ws.Range("I135").Value = 12896.375
$ws.Range("J135").Value = 158797.45
$ws.Range("K135").Value = 116067.375
$ws.Range("L135").Value = 1429177.05
$ws.Range("M135").Value = -113532.375
$ws.Range("N135").Value = -1434247.05

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 9125
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 9125
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 9125
$ws.Range("N92").Value = -12869

$ws.Range("H102").Value = 4648.857
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 4648.857
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 4648.857
$ws.Range("N102").Value = -7892.857
$ws.Range("M102").ClearContents()

$ws.Range("H113").Value = 58824616
$ws.Range("I113").Value = 90910090
$ws.Range("J113").Value = 1248.3334
$ws.Range("K113").Value = 90910090
$ws.Range("L113").Value = 1248.3334
$ws.Range("M113").Value = -90907920
$ws.Range("N113").Value = -5588.3334

$ws.Range("H122").Value = 5220727
$ws.Range("I122").Value = 3242698.5
$ws.Range("J122").Value = 25001012
$ws.Range("K122").Value = 9728095.5
$ws.Range("L122").Value = 75003036
$ws.Range("M122").Value = -9725645.5
$ws.Range("N122").Value = -75007936

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H132").Value = 4155.6855
$ws.Range("I132").Value = 7252
$ws.Range("J132").Value = 3515.0688
$ws.Range("K132").Value = 21756
$ws.Range("L132").Value = 10545.2064
$ws.Range("M132").Value = -19226

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4833426.5
$ws.Range("I22").Value = 27778954
$ws.Range("J22").Value = 2789.4736
$ws.Range("K22").Value = 27778954
$ws.Range("L22").Value = 2789.4736
$ws.Range("M22").Value = -27778659
$ws.Range("N22").Value = -3379.4736

$ws.Range("H27").Value = 4833426.5
$ws.Range("I27").Value = 27778954
$ws.Range("J27").Value = 2789.4736
$ws.Range("K27").Value = 27778954
$ws.Range("L27").Value = 2789.4736
$ws.Range("M27").Value = -27778847
$ws.Range("N27").Value = -3003.4736

$ws.Range("H93").Value = 50021692
$ws.Range("I93").Value = 24100.334
$ws.Range("J93").Value = 500000000
$ws.Range("K93").Value = 24100.334
$ws.Range("L93").Value = 500000000
$ws.Range("M93").Value = -22852.334
$ws.Range("N93").Value = -500002496

$ws.Range("H100").Value = 2122.2222
$ws.Range("I100").Value = 2253
$ws.Range("J100").Value = 2084.8572
$ws.Range("K100").Value = 2253
$ws.Range("L100").Value = 2084.8572
$ws.Range("M100").Value = -1712
$ws.Range("N100").Value = -3166.8572

$ws.Range("H127").Value = 68999
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 68999
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 68999
$ws.Range("N127").Value = -78919

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 18683.666
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 18683.666
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 18683.666
$ws.Range("N101").Value = -25173.666

$ws.Range("H103").Value = 30500
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 30500
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 30500
$ws.Range("N103").Value = -32844

$ws.Range("H105").Value = 46000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 46000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 46000
$ws.Range("N105").Value = -52988
